# OpenTbs 1.8.1 beta - rename the "xlsx*" ope keywords to the new common
# "tbs:*" keywords (shared between ODS and XLSX), and add a named cell
# example ("the_named_cell") on the "Delete me" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "Examples part 1": rename ope=xlsxNum/xlsxBool/xlsxDate -> tbs:num/tbs:bool/tbs:date
$ws1 = $wb.Worksheets.Item("Examples part 1")

$ws1.Range("E20").Value = "[a.score;ope=tbs:num]"
$ws1.Range("F20").Value = "[a.score;ope=tbs:num]"
$ws1.Range("C26").Value = "[cell2.score;block=tbs:cell;ope=tbs:num]"

$ws1.Range("C34").Value = "tbs:num"
$ws1.Range("D34").Value = "[onshow.x_num;ope=tbs:num]"
$ws1.Range("C35").Value = "tbs:bool"
$ws1.Range("D35").Value = "[onshow.x_bt;ope=tbs:bool]"
$ws1.Range("C36").Value = "tbs:date"
$ws1.Range("D36").Value = "[onshow.x_dt;ope=tbs:date]"

# Distinguish the second "Score" header (column F) from the first (column E)
$ws1.Range("F19").Value = "Score again"

# --- Sheet "Delete me": add a new example row demonstrating a named cell
$ws4 = $wb.Worksheets.Item("Delete me")
$ws4.Range("B6").Value = "And this named cell too."

# --- Workbook: define the named cell used by the new example
$wb.Names.Add("the_named_cell", "='Delete me'!`$B`$6")

Write-Host "done"
